# Update "想去人数" (want-to-go count) figures in column F across sheets.
# Values were incremented as part of a scheduled data refresh
# ("Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value  = 3037
$ws.Range("F10").Value = 7148
$ws.Range("F24").Value = 1840
$ws.Range("F31").Value = 599
$ws.Range("F32").Value = 25
$ws.Range("F36").Value = 2563
$ws.Range("F37").Value = 2849
$ws.Range("F39").Value = 61
$ws.Range("F40").Value = 199
$ws.Range("F45").Value = 349
$ws.Range("F47").Value = 199
$ws.Range("F49").Value = 67

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F25").Value = 19

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value  = 1760
$ws.Range("F8").Value  = 2825
$ws.Range("F9").Value  = 1078
$ws.Range("F10").Value = 1013
$ws.Range("F12").Value = 363
$ws.Range("F14").Value = 7802

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value  = 3037
$ws.Range("F7").Value  = 1760
$ws.Range("F9").Value  = 2825
$ws.Range("F10").Value = 7148
$ws.Range("F11").Value = 1078
$ws.Range("F12").Value = 1014
$ws.Range("F15").Value = 363
$ws.Range("F25").Value = 1840
$ws.Range("F32").Value = 599
$ws.Range("F33").Value = 25
$ws.Range("F39").Value = 2563
$ws.Range("F40").Value = 2849
$ws.Range("F42").Value = 61
$ws.Range("F43").Value = 199
$ws.Range("F45").Value = 349
$ws.Range("F46").Value = 19
$ws.Range("F47").Value = 199
